# Generate Report for Handoff
# A new source file (ffff0ac76892-4552-4a0c-bfb9-ea3d72729ca7.md) has been picked
# up and is now "Ready for handoff" alongside the existing
# b6091237-6809-4684-867b-5538749eeb17.md file (which itself was re-handed-off,
# refreshing its handoff target .xlf / timestamp). This inserts a new row on every
# sheet (pushing the ".localization-config" row down one) and refreshes the
# existing handoff metadata.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()

# Row 4 (was row 3): .localization-config / Not to be localized
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("C4").Value = "Not to be localized"

# Row 2: existing handed-off markdown file (guid changed)
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"

# Row 3 (NEW): newly discovered markdown file, also ready for handoff
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5f82ccf24aa6c800cf8aec9e86dcd0673775a23b/e2e/b6091237-6809-4684-867b-5538749eeb17.md", "", "", "b6091237-6809-4684-867b-5538749eeb17.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5f82ccf24aa6c800cf8aec9e86dcd0673775a23b/e2e/ffff0ac76892-4552-4a0c-bfb9-ea3d72729ca7.md", "", "", "ffff0ac76892-4552-4a0c-bfb9-ea3d72729ca7.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5f82ccf24aa6c800cf8aec9e86dcd0673775a23b/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()

# Row 4 (was row 3): .localization-config / Not to be localized / Ignored
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H4").Value = "Ignored"

# Row 2: refreshed handoff metadata for the existing file
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-03-01 09:38:21"
$ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H2").Value = "Include"

# Row 3 (NEW): the newly discovered file, handed off together with row 2
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-01 09:38:21"
$ws.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5f82ccf24aa6c800cf8aec9e86dcd0673775a23b/e2e/b6091237-6809-4684-867b-5538749eeb17.md", "", "", "b6091237-6809-4684-867b-5538749eeb17.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/07aa6c290f107bab9531013a8c7ff8e54395f860/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.zh-cn.xlf", "", "", "b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5f82ccf24aa6c800cf8aec9e86dcd0673775a23b/e2e/ffff0ac76892-4552-4a0c-bfb9-ea3d72729ca7.md", "", "", "ffff0ac76892-4552-4a0c-bfb9-ea3d72729ca7.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/07aa6c290f107bab9531013a8c7ff8e54395f860/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.zh-cn.xlf", "", "", "b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5f82ccf24aa6c800cf8aec9e86dcd0673775a23b/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()

# Row 4 (was row 3): .localization-config / Not to be localized / Ignored
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H4").Value = "Ignored"

# Row 2: refreshed handoff metadata for the existing file
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-03-01 09:38:31"
$ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H2").Value = "Include"

# Row 3 (NEW): the newly discovered file, handed off together with row 2
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-01 09:38:31"
$ws.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5f82ccf24aa6c800cf8aec9e86dcd0673775a23b/e2e/b6091237-6809-4684-867b-5538749eeb17.md", "", "", "b6091237-6809-4684-867b-5538749eeb17.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/127fed0d2471b9ff16a66c4222c7d99f3584afe9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.de-de.xlf", "", "", "b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5f82ccf24aa6c800cf8aec9e86dcd0673775a23b/e2e/ffff0ac76892-4552-4a0c-bfb9-ea3d72729ca7.md", "", "", "ffff0ac76892-4552-4a0c-bfb9-ea3d72729ca7.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/127fed0d2471b9ff16a66c4222c7d99f3584afe9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.de-de.xlf", "", "", "b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5f82ccf24aa6c800cf8aec9e86dcd0673775a23b/.localization-config", "", "", ".localization-config") | Out-Null

Write-Output "Report generated for handoff."
